$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# New data rows 37-41 (columns B and F hold brand-new shared strings; the
# order they are written in controls the order they land in sharedStrings.xml)
$B = @("rgdmff.tumblr.com", "lastdragonlord.tumblr.com", "offfffffffffthesouthernisles.tumblr.com", "panphangirl.tumblr.com", "typicalwelshnonsense.tumblr.com")
$F = @("Relakan", "Last Dragon", "Semangat Baru", "Terlalu Indah", "Aku Disini")
$C = @(32, 33, 32, 32, 32)
$D = @(98, 98, 98, 98, 98)
$H = @(0, 0, 29, 0, 16)
$I = @(216, 21, 255, 237, 18100)

for ($i = 0; $i -lt 5; $i++) {
    $r = 37 + $i
    $ws.Range("B$r").Value = $B[$i]
}

for ($i = 0; $i -lt 5; $i++) {
    $r = 37 + $i
    $ws.Range("F$r").Value = $F[$i]
}

for ($i = 0; $i -lt 5; $i++) {
    $r = 37 + $i
    $ws.Range("C$r").Value = $C[$i]
    $ws.Range("D$r").Value = $D[$i]
    $ws.Range("E$r").Value = "andi.sliye@yahoo.com"
    $ws.Range("H$r").Value = $H[$i]
    $ws.Range("I$r").Value = $I[$i]
}

for ($i = 0; $i -lt 5; $i++) {
    $r = 37 + $i
    # Add the mailto hyperlink on the E column, mirroring the existing rows
    $ws.Hyperlinks.Add($ws.Range("E$r"), "mailto:andi.sliye@yahoo.com")
}

# Restore the "Hyperlink" cell style used by the rest of column E (adding a
# hyperlink nudges the cell's font; put it back the way the other rows have it)
$ws.Range("E37:E41").Style = "Hyperlink"

# Rows 37, 39 and 41 in column I pick up the same visual style already used
# by several other cells in column I (copy number format/font only).
$ws.Range("I31").Copy()
$ws.Range("I37").PasteSpecial(-4122)
$ws.Range("I39").PasteSpecial(-4122)
$ws.Range("I41").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the selected cell shown when the sheet is re-opened
$ws.Range("L4").Select()
